$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC2 block currently holds "realizar a liquidação" step (B20/D20)
# TC3 block currently holds "atribuir/desatribuir" step (B28/D28)
# The commit swaps these two steps, so that TC2 now shows the
# "atribuir/desatribuir" text and TC3 shows the "realizar a liquidação" text.

$tc2Desc = $ws.Range("B20").Value2
$tc2Result = $ws.Range("D20").Value2
$tc3Desc = $ws.Range("B28").Value2
$tc3Result = $ws.Range("D28").Value2

$ws.Range("B20").Value2 = $tc3Desc
$ws.Range("D20").Value2 = $tc3Result

$ws.Range("B28").Value2 = $tc2Desc
$ws.Range("D28").Value2 = $tc2Result
